# Edit script: update EC (Estado de Cuenta) worksheet with new worker/period data.
# - Updates summary cells (Valor Mora total, Cant. Trabajadores, Cant. Periodos)
# - Replaces the 17 existing detail rows (16-32) with a new, reordered set of
#   21 detail rows (16-36): some workers removed, new worker OSNAIDER JAVIER
#   BUSTOS JARAMILLO added (6 periods), LEIDER TOMAS LORA TOBIAS gains a period.
# - Moves the two footer rows down from 37/38 to 41/42 to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room: insert 4 extra detail rows right after the last current
#    data row (32) and before the blank/footer area. Using a column-scoped
#    insert (B:J) keeps everything else (images anchored above row 6, etc.)
#    untouched, and shifts the footer rows (37/38 -> 41/42) + mergeCells
#    automatically.
# ---------------------------------------------------------------------------
$ws.Range("B32:J35").Insert(-4121)

# ---------------------------------------------------------------------------
# 2. Re-apply correct formatting:
#    - New row 36 (last detail row) must carry the special "last row" style
#      that the old row 32 used to have -> copy that format there first.
#    - Rows 32-35 (now plain/interior rows) + row 16-31 get the regular
#      interior-row format (same as row 31/16).
# ---------------------------------------------------------------------------
$ws.Range("B32:J32").Copy()
$ws.Range("B36:J36").PasteSpecial(-4122)

$ws.Range("B31:J31").Copy()
$ws.Range("B32:J35").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Write the new detail rows 16-36.
# ---------------------------------------------------------------------------
function Set-DetailRow($row, $doc, $name, $period, $mora, $salario) {
    $ws.Range("B" + $row).Value = "CC"
    $ws.Range("C" + $row).Value = $doc
    $ws.Range("D" + $row).Value = $name
    $ws.Range("E" + $row).Value = $period
    $ws.Range("F" + $row).Value = $mora
    $ws.Range("G" + $row).Value = $salario
}

Set-DetailRow 16 "73353731"    "EULISE MATTAS BARRIOS"               "2107" 8480   908526
Set-DetailRow 17 "73353731"    "EULISE MATTAS BARRIOS"               "2106" 36341  908526
Set-DetailRow 18 "73353731"    "EULISE MATTAS BARRIOS"               "2105" 30284  908526
Set-DetailRow 19 "3828458"     "LEIDER TOMAS LORA TOBIAS"             "2209" 40000  1000000
Set-DetailRow 20 "3828458"     "LEIDER TOMAS LORA TOBIAS"             "2109" 36341  1000000
Set-DetailRow 21 "3828458"     "LEIDER TOMAS LORA TOBIAS"             "2107" 8480   1000000
Set-DetailRow 22 "3828458"     "LEIDER TOMAS LORA TOBIAS"             "2106" 36341  1000000
Set-DetailRow 23 "3828458"     "LEIDER TOMAS LORA TOBIAS"             "2105" 29073  1000000
Set-DetailRow 24 "45646635"    "PAOLA MARGARITA RODRIGUEZ CABALLERO"  "1704" 7933   850000
Set-DetailRow 25 "1007323757"  "CARLOS ALFREDO RODRIGUEZ CERVANTES"   "2204" 1333   1000000
Set-DetailRow 26 "1044921366"  "OSNAIDER JAVIER BUSTOS JARAMILLO"     "2507" 36341  908526
Set-DetailRow 27 "1044921366"  "OSNAIDER JAVIER BUSTOS JARAMILLO"     "2506" 36341  908526
Set-DetailRow 28 "1044921366"  "OSNAIDER JAVIER BUSTOS JARAMILLO"     "2505" 36341  908526
Set-DetailRow 29 "1044921366"  "OSNAIDER JAVIER BUSTOS JARAMILLO"     "2504" 36341  908526
Set-DetailRow 30 "1044921366"  "OSNAIDER JAVIER BUSTOS JARAMILLO"     "2503" 36341  908526
Set-DetailRow 31 "1044921366"  "OSNAIDER JAVIER BUSTOS JARAMILLO"     "2502" 36341  908526
Set-DetailRow 32 "1044908417"  "HAWYN JESSID CASTRO SARMIENTO"        "2012" 36612  980657
Set-DetailRow 33 "1044935054"  "GUSTAVO ANDRES MARTINEZ PADILLA"      "2011" 39227  877803
Set-DetailRow 34 "1044935054"  "GUSTAVO ANDRES MARTINEZ PADILLA"      "2010" 39227  877803
Set-DetailRow 35 "1044912807"  "MIRLEYDIS DE JESUS ROMERO HURTADO"    "2101" 7268   908526
Set-DetailRow 36 "1051824856"  "CARLOS ANDRES DE ORO BUSTILLO"        "2204" 22666  1000000

# ---------------------------------------------------------------------------
# 4. Update the summary block above the table.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 597652    # VALOR MORA total
$ws.Range("C13").Value = 9         # Cant. Trabajadores
$ws.Range("F13").Value = 17        # Cant. Periodos
